$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 1.0
$ws.Range("C3").Value = 1.0
$ws.Range("D3").Value = 1.0
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 1.0

# Row 4
$ws.Range("E4").Value = 656.0

# Row 6
$ws.Range("B6").Value = 634.0

# Row 7
$ws.Range("B7").Value = 19704932.0
$ws.Range("C7").Value = 2079010.0
$ws.Range("D7").Value = 64777.0
$ws.Range("E7").Value = 7923051.0
$ws.Range("F7").Value = 353757.0

# Row 8
$ws.Range("B8").Value = 0.0
$ws.Range("C8").Value = 0.0
$ws.Range("D8").Value = 0.0
$ws.Range("E8").Value = 0.0
$ws.Range("F8").Value = 0.0

# Row 9
$ws.Range("B9").Value = 851071.0
$ws.Range("C9").Value = 3127962.0
$ws.Range("D9").Value = 114526.0
$ws.Range("E9").Value = 352860.0
$ws.Range("F9").Value = 8541412.0

# Row 10
$ws.Range("B10").Value = 0.0
$ws.Range("C10").Value = 0.0
$ws.Range("D10").Value = 0.0
$ws.Range("E10").Value = 0.0
$ws.Range("F10").Value = 0.0

# Row 11
$ws.Range("B11").Value = 232219.0
$ws.Range("C11").Value = 6672071.0
$ws.Range("D11").Value = 37758.0
$ws.Range("E11").Value = 61203.0
$ws.Range("F11").Value = 816434.0

# Row 12
$ws.Range("B12").Value = 0.0
$ws.Range("C12").Value = 0.0
$ws.Range("D12").Value = 0.0
$ws.Range("E12").Value = 0.0
$ws.Range("F12").Value = 0.0
